# Generate Report for Handback
# The "5f2c9c7b-8ee7-48fd-af0c-5866b81aa82d" handback entry (row 3) is removed
# from every sheet, and the "Correspond Handoff/Handback DateTime" stamps for
# the remaining 0e314636 entries are refreshed on the zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

function Remove-HyperlinkAt($ws, $targetAddr) {
    # Re-walk the collection from scratch and stop at the first match so
    # that deleting an item never leaves a stale/shifted iterator handle
    # pointing at the wrong hyperlink for any later items.
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $targetAddr) {
            $hl.Delete()
            return
        }
    }
}

# --- Overview sheet: drop the 5f2c9c7b row (row 3) -------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
Remove-HyperlinkAt $wsOverview '$A$3'
$wsOverview.Rows.Item(3).Delete()

# --- zh-cn sheet: refresh timestamps, drop the 5f2c9c7b row (row 3) --------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 08:47:36"
$wsZhCn.Range("H2").Value = "2016-03-22 08:47:57"
Remove-HyperlinkAt $wsZhCn '$A$3'
Remove-HyperlinkAt $wsZhCn '$D$3'
Remove-HyperlinkAt $wsZhCn '$F$3'
Remove-HyperlinkAt $wsZhCn '$G$3'
$wsZhCn.Rows.Item(3).Delete()

# --- de-de sheet: refresh timestamps, drop the 5f2c9c7b row (row 3) --------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 08:47:39"
$wsDeDe.Range("H2").Value = "2016-03-22 08:48:05"
Remove-HyperlinkAt $wsDeDe '$A$3'
Remove-HyperlinkAt $wsDeDe '$D$3'
Remove-HyperlinkAt $wsDeDe '$F$3'
Remove-HyperlinkAt $wsDeDe '$G$3'
$wsDeDe.Rows.Item(3).Delete()
